$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data to the table: Match 20, Player "Montacer", Buts 2
$ws.Range("A30").Value = 20
$ws.Range("B30").Value = "Montacer"
$ws.Range("C30").Value = 2

# Update the selection to reflect where the user clicked next (A31)
$ws.Range("A31").Select()
